$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mc3")

# Row 3 - new subject (mask size calibration data)
$ws.Range("A3").Value = 4
$ws.Range("B3").Value = "l"
$ws.Range("C3").Value = -0.797
$ws.Range("D3").Value = -0.352
$ws.Range("E3").Value = -0.877
$ws.Range("F3").Value = -0.198
$ws.Range("G3").Formula = "=AVERAGE(C3,D3)"
$ws.Range("H3").Formula = "=AVERAGE(E3,F3)"
$ws.Range("I3").Value = "e"
$ws.Range("J3").Value = "Azman"
$ws.Range("K3").Value = "  "

# Row 4 - new subject (mask size calibration data)
$ws.Range("A4").Value = 5
$ws.Range("B4").Value = "r"
$ws.Range("C4").Value = -1.057
$ws.Range("D4").Value = -0.958
$ws.Range("E4").Value = -0.997
$ws.Range("F4").Value = -1.212
$ws.Range("G4").Formula = "=AVERAGE(C4,D4)"
$ws.Range("H4").Formula = "=AVERAGE(E4,F4)"
$ws.Range("I4").Value = "e"
$ws.Range("J4").Value = "Zhi Wei"

# Move the active selection to match the author's final cursor position
[void]$ws.Range("B4").Select()
